$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1169.9166
$ws.Range("I19").Value = 940.1667
$ws.Range("K19").Value = 940.1667
$ws.Range("M19").Value = -765.1667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H46").Value = 2500
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H60").Value = 2500
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1255.6666
$ws.Range("J98").Value = 489.5
$ws.Range("L98").Value = 489.5
$ws.Range("N98").Value = -3485.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 1255.6666
$ws.Range("J122").Value = 489.5
$ws.Range("L122").Value = 1468.5
$ws.Range("N122").Value = -6368.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H126").Value = 100000
$ws.Range("J126").Value = 100000
$ws.Range("L126").Value = 100000
$ws.Range("N126").Value = -109880

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2411.2856
$ws.Range("I132").Value = 1806.4348
$ws.Range("K132").Value = 5419.3044
$ws.Range("M132").Value = -2889.3044

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 7465619
$ws.Range("I138").Value = 1333.4736
$ws.Range("J138").Value = 10420232
$ws.Range("K138").Value = 4000.4208
$ws.Range("L138").Value = 31260696
$ws.Range("M138").Value = 1139.5792
$ws.Range("N138").Value = -31270976

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H140").Value = 240390
$ws.Range("J140").Value = 240390
$ws.Range("L140").Value = 240390
$ws.Range("N140").Value = -250750

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 3074.7273
$ws.Range("I141").Value = 2992.7
$ws.Range("J141").Value = 3895
$ws.Range("K141").Value = 8978.099999999999
$ws.Range("L141").Value = 11685
$ws.Range("M141").Value = -3798.099999999999
$ws.Range("N141").Value = -22045

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4145.7393
$ws.Range("I61").Value = 3000.6086
$ws.Range("J61").Value = 5290.8696
$ws.Range("K61").Value = 3000.6086
$ws.Range("L61").Value = 5290.8696
$ws.Range("M61").Value = -2788.6086
$ws.Range("N61").Value = -5714.8696

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2858.5454
$ws.Range("I122").Value = 2806.125
$ws.Range("K122").Value = 8418.375
$ws.Range("M122").Value = -5968.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 4672.9287
$ws.Range("I132").Value = 4642.696
$ws.Range("K132").Value = 13928.088
$ws.Range("M132").Value = -11398.088

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 4145.7393
$ws.Range("I136").Value = 3000.6086
$ws.Range("J136").Value = 5290.8696
$ws.Range("K136").Value = 9001.825800000001
$ws.Range("L136").Value = 15872.6088
$ws.Range("M136").Value = -6451.825800000001
$ws.Range("N136").Value = -20972.6088

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H139").Value = 91642.28999999999
$ws.Range("J139").Value = 91642.28999999999
$ws.Range("L139").Value = 91642.28999999999
$ws.Range("N139").Value = -101922.29

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 14891.375
$ws.Range("I107").Value = 3188.6667
$ws.Range("K107").Value = 3188.6667
$ws.Range("M107").Value = -1268.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 86225.914
$ws.Range("I31").Value = 145432.42
$ws.Range("J31").Value = 3336.8
$ws.Range("K31").Value = 145432.42
$ws.Range("L31").Value = 3336.8
$ws.Range("M31").Value = -145137.42
$ws.Range("N31").Value = -3926.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 86225.914
$ws.Range("I34").Value = 145432.42
$ws.Range("J34").Value = 3336.8
$ws.Range("K34").Value = 145432.42
$ws.Range("L34").Value = 3336.8
$ws.Range("M34").Value = -145230.42
$ws.Range("N34").Value = -3740.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3852.4167
$ws.Range("I132").Value = 3929.9092
$ws.Range("K132").Value = 11789.7276
$ws.Range("M132").Value = -9259.7276

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 376.3846
$ws.Range("I86").Value = 232
$ws.Range("J86").Value = 466.625
$ws.Range("K86").Value = 696
$ws.Range("L86").Value = 1399.875
$ws.Range("M86").Value = 490
$ws.Range("N86").Value = -3771.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 21815.084
$ws.Range("I87").Value = 17722.75
$ws.Range("K87").Value = 53168.25
$ws.Range("M87").Value = -51920.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H89").Value = 376.3846
$ws.Range("I89").Value = 232
$ws.Range("J89").Value = 466.625
$ws.Range("K89").Value = 2088
$ws.Range("L89").Value = 4199.625
$ws.Range("M89").Value = 3840
$ws.Range("N89").Value = -16055.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H90").Value = 21815.084
$ws.Range("I90").Value = 17722.75
$ws.Range("K90").Value = 159504.75
$ws.Range("M90").Value = -153264.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 782.63635
$ws.Range("J113").Value = 799
$ws.Range("L113").Value = 2397
$ws.Range("N113").Value = -6737

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 40208.23
$ws.Range("I131").Value = 167849.83
$ws.Range("K131").Value = 503549.49
$ws.Range("M131").Value = -498509.49

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 57831.832
$ws.Range("I113").Value = 57831.832
$ws.Range("K113").Value = 57831.832
$ws.Range("M113").Value = -55661.832

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4366.143
$ws.Range("J122").Value = 4424.5
$ws.Range("L122").Value = 13273.5
$ws.Range("N122").Value = -18173.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2874.1538
$ws.Range("I132").Value = 2938.7083
$ws.Range("K132").Value = 8816.124899999999
$ws.Range("M132").Value = -6286.124899999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H136").Value = 30749
$ws.Range("J136").Value = 30749
$ws.Range("L136").Value = 92247
$ws.Range("N136").Value = -97347

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 7226.1816
$ws.Range("J132").Value = 7927
$ws.Range("L132").Value = 23781
$ws.Range("N132").Value = -28841

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 22003.6
$ws.Range("J101").Value = 22003.6
$ws.Range("L101").Value = 22003.6
$ws.Range("N101").Value = -28493.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3020.1667
$ws.Range("I132").Value = 2624.2
$ws.Range("K132").Value = 7872.599999999999
$ws.Range("M132").Value = -5342.599999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H141").Value = 131499.75
$ws.Range("J141").Value = 131499.75
$ws.Range("L141").Value = 131499.75
$ws.Range("N141").Value = -141859.75
